# Daily attendance processing - 2025-10-23 07:21:23
# Reorders the "Recorded By" list in column G for a set of rows so that the
# real user / backup identifier appears before the literal "System" marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "backup@backdoor.com, system, System"
    5   = "backup@backdoor.com, System"
    7   = "admin@admin.com, System"
    8   = "backup@backdoor.com, System"
    11  = "dnasr281@gmail.com, System"
    17  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, system, System"
    32  = "backup@backdoor.com, System"
    34  = "admin@admin.com, System"
    35  = "backup@backdoor.com, System"
    38  = "dnasr281@gmail.com, System"
    44  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, system, System"
    59  = "backup@backdoor.com, System"
    61  = "admin@admin.com, System"
    62  = "backup@backdoor.com, System"
    65  = "dnasr281@gmail.com, System"
    71  = "dnasr281@gmail.com, System"
    83  = "backup@backdoor.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    96  = "dnasr281@gmail.com, System"
    97  = "dnasr281@gmail.com, System"
    99  = "dnasr281@gmail.com, System"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    122 = "dnasr281@gmail.com, System"
    123 = "dnasr281@gmail.com, System"
    125 = "dnasr281@gmail.com, System"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    148 = "dnasr281@gmail.com, System"
    149 = "dnasr281@gmail.com, System"
    151 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
